$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'26.915.79"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +4.37%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.883.37"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +3.66%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'278.87"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.64%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.07%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.5372"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +4.65%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3464"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -1.51%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.06993"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +4.90%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'20.20"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +0.94%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.8100"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.54%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.07735"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -1.71%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'1.879.25"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +3.43%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'91.05"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +3.97%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'5.193"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +2.25%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'14.61"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +3.33%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.9991"
$ws.Range("D17").Style = $style

$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'0.000008061"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.46%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.9996"
$ws.Range("D19").Style = $style

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'26.974.54"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +4.30%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'2.115.37"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +3.37%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.763"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.81%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'10.09"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +0.75%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'6.218"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +2.14%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.388"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +8.51%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'147.24"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +4.14%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'17.41"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +2.04%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'1.665"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.44%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'113.84"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +3.99%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'4.378"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +0.58%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'4.338"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +2.44%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'0.08924"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +1.18%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.04944"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +1.56%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'1.183"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +4.12%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.7363"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.79%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'2.885"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +0.28%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'3.300"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +4.67%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'2.384"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +0.76%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.01861"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +0.63%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.5184"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.98%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.9616"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +0.68%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'115.98"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +3.92%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'6.206"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +0.34%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'8.182"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +1.75%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.9991"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.11%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.4514"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -1.47%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.1350"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.90%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'9.414"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +2.01%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'36.35"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("E50").Value = "  +2.10%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'1.504"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.04%  "

